$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 21.56901359558105
$ws.Range("C3").Value = 20.48993110656738
$ws.Range("C4").Value = 22.85623550415039
$ws.Range("C5").Value = 19.31190490722656
$ws.Range("C6").Value = 20.57194709777832
